# Content refresh (2025-08): new publication host + republish date, and the
# column-width relayout that Excel's "best fit" produced afterwards on the
# Elements sheet.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/measurement-context"
$wsMeta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

$wsElem = $wb.Worksheets.Item("Elements")
# Same shared string as Metadata!B2 ("URL"), so it must track the new host too.
$wsElem.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/measurement-context"
$wsElem.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/measurement-context-vs"

# Resulting best-fit column widths (stored XML width = ColumnWidth + 5/6).
$widths = @{
    1  = 15.584635416666666
    2  = 15.584635416666666
    3  = 8.959635416666666
    4  = 6.213541666666667
    5  = 4.467447916666667
    6  = 3.1197916666666665
    7  = 3.4322916666666665
    8  = 11.854166666666666
    9  = 9.678385416666666
    10 = 19.869791666666668
    11 = 13.541666666666666
    12 = 99.86979166666667
    13 = 99.86979166666667
    14 = 99.86979166666667
    15 = 11.428385416666666
    16 = 19.869791666666668
    17 = 19.869791666666668
    18 = 19.869791666666668
    19 = 19.869791666666668
    20 = 6.967447916666667
    21 = 12.776041666666666
    22 = 13.084635416666666
    23 = 14.178385416666666
    24 = 13.795572916666666
    25 = 16.248697916666668
    26 = 60.686197916666664
    27 = 4.240885416666667
    28 = 17.147135416666668
    29 = 33.744791666666664
    30 = 12.709635416666666
    31 = 10.486979166666666
    32 = 14.213541666666666
    33 = 7.389322916666667
    34 = 7.697916666666667
    35 = 99.86979166666667
    37 = 18.729166666666668
}

foreach ($col in $widths.Keys) {
    $wsElem.Columns.Item($col).ColumnWidth = $widths[$col]
}
